# Applies the "stimuli" worksheet update:
#  - adds a "generic" value in column J for the practice rows (2-5)
#  - appends a new "stim details" block (rows 27-36) describing
#    month/word_type/need_audio/need_image/word/count/find-images data

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Practice rows (2-5): new column J = "generic" ---
$ws.Range("J2").Value = "generic"
$ws.Range("J3").Value = "generic"
$ws.Range("J4").Value = "generic"
$ws.Range("J5").Value = "generic"

# --- New "stim details" block ---
$ws.Range("A27").Value = "stim details"

$ws.Range("A28").Value = "month"
$ws.Range("B28").Value = "word_type"
$ws.Range("C28").Value = "need_audio"
$ws.Range("D28").Value = "need_image"
$ws.Range("E28").Value = "word"
$ws.Range("F28").Value = "count"
$ws.Range("G28").Value = "find images"

$stimRows = @(
    @{ Row = 29; Month = 6; Type = "video" },
    @{ Row = 30; Month = 6; Type = "video" },
    @{ Row = 31; Month = 7; Type = "video" },
    @{ Row = 32; Month = 7; Type = "video" },
    @{ Row = 33; Month = 6; Type = "audio" },
    @{ Row = 34; Month = 6; Type = "audio" },
    @{ Row = 35; Month = 7; Type = "audio" },
    @{ Row = 36; Month = 7; Type = "audio" }
)

foreach ($r in $stimRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Month
    $ws.Cells.Item($r.Row, 2).Value = $r.Type
}
